$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8372160792350769
$ws.Range("B1").Value = 1.213238000869751
$ws.Range("C1").Value = 2.371771574020386
$ws.Range("D1").Value = 3.837385416030884
$ws.Range("E1").Value = 1.910194039344788
